# Commit: "made chagnes for month append"
# Appends a new month (May-24) row to the financial analysis sheet and
# updates a handful of Apr-24 figures that were recalculated as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sheet view: make rightToLeft explicit (matches rightToLeft="0") ---
$ws.DisplayRightToLeft = $false

# --- Updated figures for the existing Apr-24 row (row 2) ---
$ws.Range("B2").Value = 22021443
$ws.Range("D2").Value = 3698350
$ws.Range("G2").Value = 6321474
$ws.Range("H2").Value = "'37%"
$ws.Range("K2").Value = 13253231.02
$ws.Range("O2").Value = 3591209

# --- New row 3: May-24 ---
$ws.Range("A3").Value = "'May-24"
$ws.Range("B3").Value = 26470108
$ws.Range("C3").Value = 1242959
$ws.Range("D3").Value = 5235425
$ws.Range("E3").Value = 32948492
$ws.Range("F3").Value = 41351
$ws.Range("G3").Value = 6715975.640000001
$ws.Range("H3").Value = "'28%"
$ws.Range("I3").Value = 4649643.27
$ws.Range("J3").Value = 7875384
$ws.Range("K3").Value = 19241002.91
$ws.Range("L3").Value = "'79%"
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 261104
$ws.Range("O3").Value = 4899771
$ws.Range("P3").Value = 24401877.91
$ws.Range("Q3").Value = "'100%"
$ws.Range("R3").Value = 744437
$ws.Range("S3").Value = "'3%"
$ws.Range("T3").Value = 564180
$ws.Range("U3").Value = "'2%"
$ws.Range("V3").Value = 0
$ws.Range("W3").Value = 1308617
$ws.Range("X3").Value = "'5%"
$ws.Range("Y3").Value = 7279348.09
$ws.Range("Z3").Value = 0
$ws.Range("AA3").Value = "'0%"
$ws.Range("AB3").Value = 0
$ws.Range("AC3").Value = "'0%"
$ws.Range("AD3").Value = 7279348.09
